$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 13

# Columns A-D are text-like ("Date", "Time", "Weekday", "Week") in this sheet.
# Force them to be stored as literal text (matching existing rows) rather
# than letting Excel auto-convert date/time-looking or numeric-looking
# strings into real dates/numbers. NumberFormat is reset with ClearFormats
# afterwards so the new cells keep the sheet's default (unstyled) look,
# exactly like the other data rows.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2024-01-03"
$ws.Cells.Item($row, 1).ClearFormats()

$ws.Cells.Item($row, 2).Value = "19:03:52"
$ws.Cells.Item($row, 3).Value = "Wednesday"

$ws.Cells.Item($row, 4).NumberFormat = "@"
$ws.Cells.Item($row, 4).Value = "00"
$ws.Cells.Item($row, 4).ClearFormats()

$ws.Cells.Item($row, 5).Value = 140144
$ws.Cells.Item($row, 6).Value = 142770
$ws.Cells.Item($row, 7).Value = 172230
$ws.Cells.Item($row, 8).Value = 146977
$ws.Cells.Item($row, 9).Value = -1
$ws.Cells.Item($row, 10).Value = 117496
$ws.Cells.Item($row, 11).Value = 224012
$ws.Cells.Item($row, 12).Value = 248087
$ws.Cells.Item($row, 13).Value = 183976
$ws.Cells.Item($row, 14).Value = 109860
$ws.Cells.Item($row, 15).Value = 40168
$ws.Cells.Item($row, 16).Value = 30823
$ws.Cells.Item($row, 17).Value = 72205
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 41579
$ws.Cells.Item($row, 20).Value = -1
